# chore: update Sheets via scheduled runner
# Refresh currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ
# (columns H/I/J) and the dependent Leve price/profit columns (K/L/M/N) on
# each crafting-sheet with newly pulled market data.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 426.26666
$ws.Range("J38").Value = 2799.5
$ws.Range("L38").Value = 8398.5
$ws.Range("N38").Value = -9142.5
$ws.Range("H100").Value = 2053.2307
$ws.Range("I100").Value = 1704.5349
$ws.Range("J100").Value = 3719.2222
$ws.Range("K100").Value = 1704.5349
$ws.Range("L100").Value = 3719.2222
$ws.Range("M100").Value = -1163.5349
$ws.Range("N100").Value = -4801.2222
$ws.Range("H116").Value = 603600
$ws.Range("I116").Value = 752750
$ws.Range("K116").Value = 752750
$ws.Range("M116").Value = -749308
$ws.Range("H132").Value = 3304.2903
$ws.Range("I132").Value = 2459.6611
$ws.Range("K132").Value = 7378.9833
$ws.Range("M132").Value = -4848.9833
$ws.Range("H137").Value = 5827.294
$ws.Range("I137").Value = 2208.9
$ws.Range("K137").Value = 6626.700000000001
$ws.Range("M137").Value = -4076.700000000001
$ws.Range("H138").Value = 6376.387
$ws.Range("I138").Value = 2552.1667
$ws.Range("J138").Value = 7294.2
$ws.Range("K138").Value = 7656.500100000001
$ws.Range("L138").Value = 21882.6
$ws.Range("M138").Value = -2516.500100000001
$ws.Range("N138").Value = -32162.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 463.4375
$ws.Range("I5").Value = 472.08334
$ws.Range("J5").Value = 437.5
$ws.Range("K5").Value = 472.08334
$ws.Range("L5").Value = 437.5
$ws.Range("M5").Value = -360.08334
$ws.Range("N5").Value = -661.5
$ws.Range("H102").Value = 3811.4211
$ws.Range("I102").Value = 3800.9443
$ws.Range("K102").Value = 3800.9443
$ws.Range("M102").Value = -2178.9443
$ws.Range("H132").Value = 21209.436
$ws.Range("I132").Value = 23161.734
$ws.Range("J132").Value = 5265.6665
$ws.Range("K132").Value = 69485.202
$ws.Range("L132").Value = 15796.9995
$ws.Range("M132").Value = -66955.202
$ws.Range("N132").Value = -20856.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 463.4375
$ws.Range("I4").Value = 472.08334
$ws.Range("J4").Value = 437.5
$ws.Range("K4").Value = 472.08334
$ws.Range("L4").Value = 437.5
$ws.Range("M4").Value = -357.08334
$ws.Range("N4").Value = -667.5
$ws.Range("H107").Value = 2600.25
$ws.Range("I107").Value = 2334.1667
$ws.Range("J107").Value = 3398.5
$ws.Range("K107").Value = 2334.1667
$ws.Range("L107").Value = 3398.5
$ws.Range("M107").Value = -414.1667000000002
$ws.Range("N107").Value = -7238.5
$ws.Range("H134").Value = 2085.9167
$ws.Range("I134").Value = 1920.8392
$ws.Range("K134").Value = 5762.517599999999
$ws.Range("M134").Value = -3227.517599999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1044.7693
$ws.Range("I16").Value = 871.1818
$ws.Range("K16").Value = 871.1818
$ws.Range("M16").Value = -584.1818
$ws.Range("H31").Value = 4994.625
$ws.Range("I31").Value = 3653.9524
$ws.Range("K31").Value = 3653.9524
$ws.Range("M31").Value = -3358.9524
$ws.Range("H34").Value = 4994.625
$ws.Range("I34").Value = 3653.9524
$ws.Range("K34").Value = 3653.9524
$ws.Range("M34").Value = -3451.9524
$ws.Range("H58").Value = 58530.555
$ws.Range("I58").Value = 69856.664
$ws.Range("J58").Value = 1900
$ws.Range("K58").Value = 69856.664
$ws.Range("L58").Value = 1900
$ws.Range("M58").Value = -69653.664
$ws.Range("N58").Value = -2306
$ws.Range("H86").Value = 29049.176
$ws.Range("I86").Value = 34010.285
$ws.Range("K86").Value = 34010.285
$ws.Range("M86").Value = -32887.285
$ws.Range("H89").Value = 29049.176
$ws.Range("I89").Value = 34010.285
$ws.Range("K89").Value = 170051.425
$ws.Range("M89").Value = -164435.425
$ws.Range("H107").Value = 143256.58
$ws.Range("I107").Value = 167066
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 167066
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = -165146
$ws.Range("N107").Value = -4240
$ws.Range("H113").Value = 1044.7693
$ws.Range("I113").Value = 871.1818
$ws.Range("K113").Value = 871.1818
$ws.Range("M113").Value = 1298.8182
$ws.Range("H132").Value = 2776.5454
$ws.Range("I132").Value = 2655.3
$ws.Range("K132").Value = 7965.900000000001
$ws.Range("M132").Value = -5435.900000000001
$ws.Range("H134").Value = 32231.03
$ws.Range("I134").Value = 38447.07
$ws.Range("K134").Value = 115341.21
$ws.Range("M134").Value = -112806.21
$ws.Range("H136").Value = 58530.555
$ws.Range("I136").Value = 69856.664
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 209569.992
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -207019.992
$ws.Range("N136").Value = -10800
$ws.Range("H141").Value = 392573
$ws.Range("J141").Value = 392573
$ws.Range("L141").Value = 392573
$ws.Range("N141").Value = -402933

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 94.80768999999999
$ws.Range("I2").Value = 33.75
$ws.Range("K2").Value = 202.5
$ws.Range("M2").Value = -89.5
$ws.Range("H5").Value = 2123.9285
$ws.Range("I5").Value = 944.75
$ws.Range("J5").Value = 3008.3125
$ws.Range("K5").Value = 2834.25
$ws.Range("L5").Value = 9024.9375
$ws.Range("M5").Value = -2722.25
$ws.Range("N5").Value = -9248.9375
$ws.Range("H12").Value = 9.800000000000001
$ws.Range("I12").Value = 2.3333333
$ws.Range("J12").Value = 13
$ws.Range("K12").Value = 6.999999900000001
$ws.Range("L12").Value = 39
$ws.Range("M12").Value = 166.0000001
$ws.Range("N12").Value = -385
$ws.Range("H80").Value = 2994.75
$ws.Range("I80").Value = 2621.5
$ws.Range("J80").Value = 3069.4
$ws.Range("K80").Value = 7864.5
$ws.Range("L80").Value = 9208.200000000001
$ws.Range("M80").Value = -6928.5
$ws.Range("N80").Value = -11080.2
$ws.Range("H83").Value = 2994.75
$ws.Range("I83").Value = 2621.5
$ws.Range("J83").Value = 3069.4
$ws.Range("K83").Value = 23593.5
$ws.Range("L83").Value = 27624.6
$ws.Range("M83").Value = -18913.5
$ws.Range("N83").Value = -36984.60000000001
$ws.Range("H135").Value = 2123.9285
$ws.Range("I135").Value = 944.75
$ws.Range("J135").Value = 3008.3125
$ws.Range("K135").Value = 8502.75
$ws.Range("L135").Value = 27074.8125
$ws.Range("M135").Value = -5967.75
$ws.Range("N135").Value = -32144.8125

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3469.077
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 4085.4285
$ws.Range("K80").Value = 2750
$ws.Range("L80").Value = 4085.4285
$ws.Range("M80").Value = -1752
$ws.Range("N80").Value = -6081.4285
$ws.Range("H83").Value = 3469.077
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 4085.4285
$ws.Range("K83").Value = 13750
$ws.Range("L83").Value = 20427.1425
$ws.Range("M83").Value = -8758
$ws.Range("N83").Value = -30411.1425
$ws.Range("H122").Value = 3289.2104
$ws.Range("I122").Value = 1975.7059
$ws.Range("J122").Value = 4352.524
$ws.Range("K122").Value = 5927.1177
$ws.Range("L122").Value = 13057.572
$ws.Range("M122").Value = -3477.1177
$ws.Range("N122").Value = -17957.572

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 353260.44
$ws.Range("I7").Value = 378983.44
$ws.Range("K7").Value = 378983.44
$ws.Range("M7").Value = -378871.44
$ws.Range("H126").Value = 353260.44
$ws.Range("I126").Value = 378983.44
$ws.Range("K126").Value = 1136950.32
$ws.Range("M126").Value = -1134480.32
$ws.Range("H132").Value = 43485.633
$ws.Range("I132").Value = 57117.047
$ws.Range("J132").Value = 5999.25
$ws.Range("K132").Value = 171351.141
$ws.Range("L132").Value = 17997.75
$ws.Range("M132").Value = -168821.141
$ws.Range("N132").Value = -23057.75
$ws.Range("H137").Value = 79999
$ws.Range("J137").Value = 79999
$ws.Range("L137").Value = 79999
$ws.Range("N137").Value = -90199

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 49635.19
$ws.Range("I126").Value = 49635.19
$ws.Range("K126").Value = 148905.57
$ws.Range("M126").Value = -146435.57
$ws.Range("H132").Value = 35608.777
$ws.Range("I132").Value = 42945.85
$ws.Range("K132").Value = 128837.55
$ws.Range("M132").Value = -126307.55
